$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells to match new custom-method column names
$ws.Range("C1").Value = "ActionKeyword"
$ws.Range("D1").Value = "ElementName"

# Update the active selection to D4
$ws.Range("D4").Select()
